# Commit: "Added User Story descrip docx"
#
# Inserts a new row into the "MyBar Artifacts" table, right above the
# "User Guide / Manual" entry (old row 81), recording a new artifact:
#   General/Particular Context : ">"  (same bullet marker used throughout)
#   Particular Context         : "User Stories Description"
#   Main Contributor           : "Adam Clark"
# Every row below shifts down by one, the Table1 ListObject grows by one
# row, and the view selection moves to the cell that was being edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row above the current row 81 ("User Guide / Manual"),
#        pushing it (and everything after it) down to row 82.
$ws.Rows.Item(81).Insert()
$ws.Rows.Item(81).RowHeight = 15.75

# --- 2. Populate the new row's cells.
$ws.Range("A81").Value = ">"
$ws.Range("B81").Value = "User Stories Description"
$ws.Range("C81").Value = "Adam Clark"

# --- 3. Grow the worksheet table (Table1) by one row so the new record
#        is included (was A2:F111, now A2:F112).
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:F112"))

# --- 4. Leave the selection on the cell that was being worked on
#        (mirrors the saved view state after the edit).
$ws.Range("C82").Select()
